$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "25.950.76"
$ws.Range("E2").Value = "  -1.14%  "

# Row 3 - Ethereum
Set-TextCell "D3" "1.638.28"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.24%  "

# Row 5 - BNB
Set-TextCell "D5" "215.69"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.18%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.32%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.13%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -1.74%  "

# Row 11 - TRON
Set-TextCell "D11" "0.0795"
$ws.Range("E11").Value = "  +0.07%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "1.866.34"

# Row 13 - Polkadot
Set-TextCell "D13" "4.27"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14 - WrappedEther
Set-TextCell "D14" "1.643.12"
$ws.Range("E14").Value = "  -0.82%  "

# Row 15 - Polygon
Set-TextCell "D15" "0.543"
$ws.Range("E15").Value = "  -0.60%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.14%  "

# Row 17 - Litecoin
Set-TextCell "D17" "62.91"
$ws.Range("E17").Value = "  -0.52%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "25.936.78"
$ws.Range("E18").Value = "  -1.15%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.31%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "192.67"
$ws.Range("E20").Value = "  -1.24%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.86%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -1.30%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  -0.71%  "

# Row 24 - Stellar
Set-TextCell "D24" "0.131"
$ws.Range("E24").Value = "  +5.31%  "

# Row 25 - Toncoin
Set-TextCell "D25" "1.78"
$ws.Range("E25").Value = "  -0.52%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.20%  "

# Row 27 - Monero
Set-TextCell "D27" "143.23"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28 - Cosmos
Set-TextCell "D28" "6.88"
$ws.Range("E28").Value = "  -1.25%  "

# Row 29 - EthereumClassic
Set-TextCell "D29" "15.57"
$ws.Range("E29").Value = "  -0.15%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.54%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -1.44%  "

# Row 33 - Filecoin
Set-TextCell "D33" "3.25"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -4.78%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.84%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  -1.05%  "

# Row 37 - Maker
Set-TextCell "D37" "1.132.49"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.70%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.63%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.20%  "

# Row 41 - now FraxShare (was PaxDollar)
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D41" "5.48"
$ws.Range("E41").Value = "  -0.88%  "

# Row 42 - now Quant (was FraxShare)
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D42" "99.30"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43 - now TrustWalletToken (was Quant)
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D43" "0.797"
$ws.Range("E43").Value = "  -0.31%  "

# Row 44 - now RocketPoolETH (was TrustWalletToken)
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D44" "1.776.06"
$ws.Range("E44").Value = "  -0.36%  "

# Row 45 - now BabyDogeCoin (was RocketPoolETH)
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D45" "0.0₆0115"
$ws.Range("E45").Value = "  +2.29%  "

# Row 46 - now Aave (was BabyDogeCoin)
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D46" "56.61"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47 - now Cronos (was Aave)
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D47" "0.0530"
$ws.Range("E47").Value = "  +2.21%  "

# Row 48 - now RenderToken (was Cronos)
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D48" "1.47"
$ws.Range("E48").Value = "  -1.08%  "

# Row 49 - now EnergySwap (was RenderToken)
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "7.69"
$ws.Range("E49").Value = "  +0.30%  "

# Row 50 - now Mantle (was EnergySwap)
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D50" "0.415"
$ws.Range("E50").Value = "  -0.75%  "

# Row 51 - now Algorand (was Mantle)
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D51" "0.0960"
$ws.Range("E51").Value = "  -0.79%  "
